{"js": "const body = context.document.body;\nconst results = body.search(\"interview_intro_prompt\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"interview.intro_prompt\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"interview_intro_prompt\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"interview.intro_prompt\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
